$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Artn"
$ws.Range("C2").Value = "Ret"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.306751666666667
$ws.Range("H2").Value = 6.920255
$ws.Range("I2").Value = 0.8617934705859058
$ws.Range("J2").Value = 0.8617934705859057
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.273058
$ws.Range("N2").Value = 3.819174
$ws.Range("O2").Value = 0.1826308343983922
$ws.Range("P2").Value = 0.1826308343983922
$ws.Range("Q2").Value = 2.936628663263333
$ws.Range("R2").Value = 26.42965796937
$ws.Range("S2").Value = 0.1573900606121903
$ws.Range("T2").Value = 0.1573900606121902

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Artn"
$ws.Range("C3").Value = "Ret"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.306751666666667
$ws.Range("H3").Value = 6.920255
$ws.Range("I3").Value = 0.8617934705859058
$ws.Range("J3").Value = 0.8617934705859057
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.094146333333334
$ws.Range("N3").Value = 15.282439
$ws.Range("O3").Value = 0.7307979647464429
$ws.Range("P3").Value = 0.7307979647464428
$ws.Range("Q3").Value = 11.75093054466056
$ws.Range("R3").Value = 105.758374901945
$ws.Range("S3").Value = 0.6297969143359534
$ws.Range("T3").Value = 0.6297969143359533

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Artn"
$ws.Range("C4").Value = "Ret"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.306751666666667
$ws.Range("H4").Value = 6.920255
$ws.Range("I4").Value = 0.8617934705859058
$ws.Range("J4").Value = 0.8617934705859057
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.012811
$ws.Range("N4").Value = 0.038433
$ws.Range("O4").Value = 0.00183784526665541
$ws.Range("P4").Value = 0.00183784526665541
$ws.Range("Q4").Value = 0.02955179560166667
$ws.Range("R4").Value = 0.265966160415
$ws.Range("S4").Value = 0.001583843050750845
$ws.Range("T4").Value = 0.001583843050750845

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Artn"
$ws.Range("C5").Value = "Ret"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.306751666666667
$ws.Range("H5").Value = 6.920255
$ws.Range("I5").Value = 0.8617934705859058
$ws.Range("J5").Value = 0.8617934705859057
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5906476666666668
$ws.Range("N5").Value = 1.771943
$ws.Range("O5").Value = 0.08473335558850956
$ws.Range("P5").Value = 0.08473335558850954
$ws.Range("Q5").Value = 1.362477489496111
$ws.Range("R5").Value = 12.262297405465
$ws.Range("S5").Value = 0.07302265258701131
$ws.Range("T5").Value = 0.07302265258701129

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Artn"
$ws.Range("C6").Value = "Ret"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3699356666666667
$ws.Range("H6").Value = 1.109807
$ws.Range("I6").Value = 0.1382065294140942
$ws.Range("J6").Value = 0.1382065294140942
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.273058
$ws.Range("N6").Value = 3.819174
$ws.Range("O6").Value = 0.1826308343983922
$ws.Range("P6").Value = 0.1826308343983922
$ws.Range("Q6").Value = 0.4709495599353333
$ws.Range("R6").Value = 4.238546039418
$ws.Range("S6").Value = 0.02524077378620196
$ws.Range("T6").Value = 0.02524077378620196

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Artn"
$ws.Range("C7").Value = "Ret"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3699356666666667
$ws.Range("H7").Value = 1.109807
$ws.Range("I7").Value = 0.1382065294140942
$ws.Range("J7").Value = 0.1382065294140942
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.094146333333334
$ws.Range("N7").Value = 15.282439
$ws.Range("O7").Value = 0.7307979647464429
$ws.Range("P7").Value = 0.7307979647464428
$ws.Range("Q7").Value = 1.884506419919222
$ws.Range("R7").Value = 16.960557779273
$ws.Range("S7").Value = 0.1010010504104894
$ws.Range("T7").Value = 0.1010010504104894

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Artn"
$ws.Range("C8").Value = "Ret"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.3699356666666667
$ws.Range("H8").Value = 1.109807
$ws.Range("I8").Value = 0.1382065294140942
$ws.Range("J8").Value = 0.1382065294140942
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.012811
$ws.Range("N8").Value = 0.038433
$ws.Range("O8").Value = 0.00183784526665541
$ws.Range("P8").Value = 0.00183784526665541
$ws.Range("Q8").Value = 0.004739245825666667
$ws.Range("R8").Value = 0.042653212431
$ws.Range("S8").Value = 0.0002540022159045647
$ws.Range("T8").Value = 0.0002540022159045647

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Artn"
$ws.Range("C9").Value = "Ret"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.3699356666666667
$ws.Range("H9").Value = 1.109807
$ws.Range("I9").Value = 0.1382065294140942
$ws.Range("J9").Value = 0.1382065294140942
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5906476666666668
$ws.Range("N9").Value = 1.771943
$ws.Range("O9").Value = 0.08473335558850956
$ws.Range("P9").Value = 0.08473335558850954
$ws.Range("Q9").Value = 0.2185016383334445
$ws.Range("R9").Value = 1.966514745001
$ws.Range("S9").Value = 0.01171070300149825
$ws.Range("T9").Value = 0.01171070300149825
